$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 804.22
$ws.Range("J17").Value = 804.22
$ws.Range("L17").Value = 2412.66
$ws.Range("N17").Value = -2748.66
$ws.Range("H41").Value = 3632.3333
$ws.Range("I41").Value = 3164.8333
$ws.Range("J41").Value = 4099.8335
$ws.Range("K41").Value = 3164.8333
$ws.Range("L41").Value = 4099.8335
$ws.Range("M41").Value = -2724.8333
$ws.Range("N41").Value = -4979.8335
$ws.Range("H86").Value = 2023396.8
$ws.Range("I86").Value = 2696256.2
$ws.Range("K86").Value = 2696256.2
$ws.Range("M86").Value = -2695133.2
$ws.Range("H89").Value = 2023396.8
$ws.Range("I89").Value = 2696256.2
$ws.Range("K89").Value = 13481281
$ws.Range("M89").Value = -13475665
$ws.Range("H107").Value = 815.63635
$ws.Range("I107").Value = 828
$ws.Range("J107").Value = 760
$ws.Range("K107").Value = 828
$ws.Range("L107").Value = 760
$ws.Range("M107").Value = 1092
$ws.Range("N107").Value = -4600
$ws.Range("H112").Value = 4061.3103
$ws.Range("J112").Value = 4236.222
$ws.Range("L112").Value = 12708.666
$ws.Range("N112").Value = -14924.666
$ws.Range("H137").Value = 8551037
$ws.Range("I137").Value = 844
$ws.Range("K137").Value = 2532
$ws.Range("M137").Value = 18
$ws.Range("H138").Value = 4670.772
$ws.Range("I138").Value = 1083.1765
$ws.Range("J138").Value = 6195.5
$ws.Range("K138").Value = 3249.5295
$ws.Range("L138").Value = 18586.5
$ws.Range("M138").Value = 1890.4705
$ws.Range("N138").Value = -28866.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6023.391
$ws.Range("I32").Value = 3260.5945
$ws.Range("K32").Value = 3260.5945
$ws.Range("M32").Value = -2973.5945
$ws.Range("H61").Value = 11143.765
$ws.Range("I61").Value = 18607
$ws.Range("J61").Value = 2747.625
$ws.Range("K61").Value = 18607
$ws.Range("L61").Value = 2747.625
$ws.Range("M61").Value = -18395
$ws.Range("N61").Value = -3171.625
$ws.Range("H88").Value = 72118.57000000001
$ws.Range("I88").Value = 832.5
$ws.Range("K88").Value = 832.5
$ws.Range("M88").Value = -426.5
$ws.Range("H91").Value = 72118.57000000001
$ws.Range("I91").Value = 832.5
$ws.Range("K91").Value = 832.5
$ws.Range("M91").Value = 571.5
$ws.Range("H110").Value = 853092.9
$ws.Range("I110").Value = 1362254.2
$ws.Range("K110").Value = 1362254.2
$ws.Range("M110").Value = -1360209.2
$ws.Range("H122").Value = 3203.64
$ws.Range("I122").Value = 1654.6
$ws.Range("K122").Value = 4963.799999999999
$ws.Range("M122").Value = -2513.799999999999
$ws.Range("H136").Value = 11143.765
$ws.Range("I136").Value = 18607
$ws.Range("J136").Value = 2747.625
$ws.Range("K136").Value = 55821
$ws.Range("L136").Value = 8242.875
$ws.Range("M136").Value = -53271
$ws.Range("N136").Value = -13342.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1616.381
$ws.Range("I107").Value = 1584.5
$ws.Range("K107").Value = 1584.5
$ws.Range("M107").Value = 335.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 287.9091
$ws.Range("I7").Value = 185
$ws.Range("K7").Value = 185
$ws.Range("M7").Value = -72
$ws.Range("H31").Value = 5699.2085
$ws.Range("I31").Value = 1850.1364
$ws.Range("K31").Value = 1850.1364
$ws.Range("M31").Value = -1555.1364
$ws.Range("H34").Value = 5699.2085
$ws.Range("I34").Value = 1850.1364
$ws.Range("K34").Value = 1850.1364
$ws.Range("M34").Value = -1648.1364
$ws.Range("H62").Value = 63360
$ws.Range("I62").Value = 5400
$ws.Range("J62").Value = 77850
$ws.Range("K62").Value = 5400
$ws.Range("L62").Value = 77850
$ws.Range("M62").Value = -4776
$ws.Range("N62").Value = -79098
$ws.Range("H65").Value = 63360
$ws.Range("I65").Value = 5400
$ws.Range("J65").Value = 77850
$ws.Range("K65").Value = 27000
$ws.Range("L65").Value = 389250
$ws.Range("M65").Value = -23880
$ws.Range("N65").Value = -395490
$ws.Range("H86").Value = 6596.25
$ws.Range("I86").Value = 4795
$ws.Range("J86").Value = 12000
$ws.Range("K86").Value = 4795
$ws.Range("L86").Value = 12000
$ws.Range("M86").Value = -3672
$ws.Range("N86").Value = -14246
$ws.Range("H89").Value = 6596.25
$ws.Range("I89").Value = 4795
$ws.Range("J89").Value = 12000
$ws.Range("K89").Value = 23975
$ws.Range("L89").Value = 60000
$ws.Range("M89").Value = -18359
$ws.Range("N89").Value = -71232
$ws.Range("H94").Value = 2625
$ws.Range("J94").Value = 2534.5833
$ws.Range("L94").Value = 2534.5833
$ws.Range("N94").Value = -3436.5833
$ws.Range("H107").Value = 1136839.5
$ws.Range("I107").Value = 1515561.2
$ws.Range("K107").Value = 1515561.2
$ws.Range("M107").Value = -1513641.2
$ws.Range("H132").Value = 13908782
$ws.Range("I132").Value = 16679789
$ws.Range("J132").Value = 53746
$ws.Range("K132").Value = 50039367
$ws.Range("L132").Value = 161238
$ws.Range("M132").Value = -50036837
$ws.Range("N132").Value = -166298

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 44041456
$ws.Range("J4").Value = 152861570
$ws.Range("L4").Value = 458584710
$ws.Range("N4").Value = -458584934
$ws.Range("H12").Value = 133.625
$ws.Range("J12").Value = 87
$ws.Range("L12").Value = 261
$ws.Range("N12").Value = -607
$ws.Range("H107").Value = 1279.8125
$ws.Range("J107").Value = 1326.4667
$ws.Range("L107").Value = 3979.4001
$ws.Range("N107").Value = -7819.4001
$ws.Range("H123").Value = 17499.834
$ws.Range("J123").Value = 24999.5
$ws.Range("L123").Value = 74998.5
$ws.Range("N123").Value = -79898.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4335232.5
$ws.Range("I70").Value = 9529411
$ws.Range("J70").Value = 6750
$ws.Range("K70").Value = 9529411
$ws.Range("L70").Value = 6750
$ws.Range("M70").Value = -9529141
$ws.Range("N70").Value = -7290
$ws.Range("H73").Value = 4335232.5
$ws.Range("I73").Value = 9529411
$ws.Range("J73").Value = 6750
$ws.Range("K73").Value = 9529411
$ws.Range("L73").Value = 6750
$ws.Range("M73").Value = -9528475
$ws.Range("N73").Value = -8622
$ws.Range("H102").Value = 18525730
$ws.Range("I102").Value = 26322884
$ws.Range("K102").Value = 26322884
$ws.Range("M102").Value = -26321262
$ws.Range("H122").Value = 5820.125
$ws.Range("J122").Value = 6831.8335
$ws.Range("L122").Value = 20495.5005
$ws.Range("N122").Value = -25395.5005
$ws.Range("H126").Value = 4290.885
$ws.Range("I126").Value = 3102.2144
$ws.Range("J126").Value = 5677.6665
$ws.Range("K126").Value = 9306.643199999999
$ws.Range("L126").Value = 17032.9995
$ws.Range("M126").Value = -6836.643199999999
$ws.Range("N126").Value = -21972.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5246.4165
$ws.Range("I7").Value = 1984
$ws.Range("K7").Value = 1984
$ws.Range("M7").Value = -1872
$ws.Range("H22").Value = 1167.2858
$ws.Range("I22").Value = 993.25
$ws.Range("J22").Value = 1399.3334
$ws.Range("K22").Value = 993.25
$ws.Range("L22").Value = 1399.3334
$ws.Range("M22").Value = -698.25
$ws.Range("N22").Value = -1989.3334
$ws.Range("H27").Value = 1167.2858
$ws.Range("I27").Value = 993.25
$ws.Range("J27").Value = 1399.3334
$ws.Range("K27").Value = 993.25
$ws.Range("L27").Value = 1399.3334
$ws.Range("M27").Value = -886.25
$ws.Range("N27").Value = -1613.3334
$ws.Range("H61").Value = 4996.6
$ws.Range("I61").Value = 4996.5
$ws.Range("K61").Value = 4996.5
$ws.Range("M61").Value = -4794.5
$ws.Range("H93").Value = 2239.5
$ws.Range("I93").Value = 538.3
$ws.Range("J93").Value = 6492.5
$ws.Range("K93").Value = 538.3
$ws.Range("L93").Value = 6492.5
$ws.Range("M93").Value = 709.7
$ws.Range("N93").Value = -8988.5
$ws.Range("H113").Value = 4996.6
$ws.Range("I113").Value = 4996.5
$ws.Range("K113").Value = 4996.5
$ws.Range("M113").Value = -2826.5
$ws.Range("H122").Value = 54429784
$ws.Range("I122").Value = 83337120
$ws.Range("J122").Value = 15886676
$ws.Range("K122").Value = 250011360
$ws.Range("L122").Value = 47660028
$ws.Range("M122").Value = -250008910
$ws.Range("N122").Value = -47664928
$ws.Range("H126").Value = 5246.4165
$ws.Range("I126").Value = 1984
$ws.Range("K126").Value = 5952
$ws.Range("M126").Value = -3482
$ws.Range("H136").Value = 2133.95
$ws.Range("I136").Value = 1961.2572
$ws.Range("K136").Value = 5883.7716
$ws.Range("M136").Value = -3333.7716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1958.8462
$ws.Range("I107").Value = 2257.7
$ws.Range("J107").Value = 962.6667
$ws.Range("K107").Value = 6773.099999999999
$ws.Range("L107").Value = 2888.0001
$ws.Range("M107").Value = -4853.099999999999
$ws.Range("N107").Value = -6728.0001
$ws.Range("H115").Value = 79999.5
$ws.Range("J115").Value = 79999.5
$ws.Range("L115").Value = 79999.5
$ws.Range("N115").Value = -83133.5
$ws.Range("H122").Value = 5321.3794
$ws.Range("I122").Value = 4469.227
$ws.Range("J122").Value = 7999.5713
$ws.Range("K122").Value = 13407.681
$ws.Range("L122").Value = 23998.7139
$ws.Range("M122").Value = -10957.681
$ws.Range("N122").Value = -28898.7139
$ws.Range("H136").Value = 8430.803
$ws.Range("I136").Value = 3324.7
$ws.Range("K136").Value = 9974.099999999999
$ws.Range("M136").Value = -7424.099999999999
